$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Real student names that replace the placeholder "Fname Lname" shared string
$names = @(
    "George Stevenson",
    "Ann Smith",
    "Matt Smith",
    "Jennifer Mikar",
    "Jennifer Mikar",
    "Matt Johnson",
    "Ann Cramer",
    "Jennifer Stevenson",
    "Jennifer Stevenson",
    "Sam Cramer",
    "Matt Johnson",
    "Jennifer Smith",
    "George Johnson",
    "Matt Smith",
    "Ann Mikar",
    "Matt Stevenson",
    "Jennifer Smith",
    "Matt Stevenson",
    "Sam Mikar",
    "Ann Smith",
    "Matt Cramer",
    "Ann Stevenson",
    "Jennifer Johnson",
    "Sam Johnson",
    "Jennifer Stevenson",
    "Ann Cramer",
    "Matt Mikar",
    "George Johnson",
    "Jennifer Cramer",
    "Sam Smith",
    "George Johnson",
    "Ann Johnson",
    "George Johnson",
    "Matt Johnson",
    "Matt Cramer",
    "George Stevenson",
    "Sam Johnson",
    "Matt Stevenson",
    "Jennifer Stevenson",
    "George Cramer",
    "Matt Stevenson",
    "Ann Smith",
    "Ann Stevenson",
    "Matt Cramer",
    "George Mikar",
    "Matt Smith",
    "Matt Stevenson",
    "Matt Smith",
    "Ann Mikar",
    "Matt Stevenson",
    "Sam Johnson",
    "Matt Johnson",
    "Jennifer Smith",
    "Jennifer Johnson",
    "Ann Stevenson",
    "George Smith",
    "Ann Johnson",
    "Jennifer Johnson",
    "George Johnson",
    "Matt Johnson",
    "Ann Smith",
    "Matt Mikar",
    "Matt Johnson",
    "George Smith",
    "Sam Smith",
    "Ann Stevenson",
    "Ann Smith",
    "Jennifer Smith",
    "Jennifer Cramer",
    "Jennifer Mikar",
    "Ann Mikar",
    "George Stevenson",
    "Sam Stevenson",
    "George Mikar",
    "Jennifer Johnson",
    "Jennifer Johnson",
    "Sam Johnson",
    "George Johnson",
    "Ann Stevenson",
    "Sam Johnson",
    "Matt Johnson",
    "Sam Johnson",
    "Jennifer Mikar",
    "Jennifer Mikar",
    "Ann Stevenson",
    "Sam Smith",
    "Matt Johnson",
    "Matt Cramer",
    "Sam Mikar",
    "Ann Stevenson",
    "Matt Stevenson",
    "George Cramer",
    "Sam Johnson",
    "Jennifer Mikar",
    "Ann Stevenson",
    "Jennifer Mikar",
    "Sam Stevenson",
    "Sam Cramer",
    "Sam Mikar",
    "Sam Johnson",
    "Ann Stevenson",
    "Ann Stevenson",
    "Jennifer Stevenson",
    "Jennifer Stevenson",
    "Jennifer Smith",
    "Matt Johnson",
    "Sam Cramer",
    "George Mikar",
    "Sam Johnson",
    "Jennifer Stevenson",
    "Matt Johnson",
    "Jennifer Smith",
    "George Smith",
    "Ann Smith",
    "Ann Stevenson",
    "Jennifer Smith",
    "George Johnson",
    "Ann Mikar",
    "Jennifer Smith",
    "Matt Smith",
    "Sam Smith",
    "Jennifer Cramer",
    "George Mikar",
    "Sam Mikar",
    "Sam Cramer",
    "George Mikar",
    "Matt Mikar",
    "Sam Stevenson",
    "Sam Cramer",
    "Sam Mikar",
    "Jennifer Mikar",
    "Ann Johnson",
    "Jennifer Cramer",
    "George Johnson",
    "George Mikar",
    "Sam Mikar",
    "Matt Johnson",
    "Ann Smith",
    "Matt Stevenson",
    "Matt Johnson",
    "Sam Mikar",
    "Ann Cramer",
    "Sam Johnson",
    "Jennifer Mikar",
    "Matt Johnson",
    "Ann Stevenson",
    "Jennifer Johnson",
    "Jennifer Stevenson",
    "Sam Cramer",
    "Sam Stevenson",
    "Ann Mikar",
    "Matt Stevenson",
    "Ann Johnson",
    "Jennifer Mikar",
    "Sam Johnson",
    "Sam Smith",
    "Ann Mikar",
    "Matt Stevenson",
    "Matt Cramer",
    "Jennifer Mikar",
    "George Johnson",
    "Sam Mikar",
    "Sam Cramer",
    "Sam Johnson",
    "Matt Smith",
    "Matt Stevenson",
    "Sam Stevenson",
    "Jennifer Mikar",
    "Jennifer Mikar",
    "George Smith",
    "Matt Mikar",
    "Sam Johnson",
    "George Stevenson",
    "Matt Stevenson",
    "Matt Smith",
    "Sam Smith",
    "Sam Smith",
    "Matt Stevenson",
    "Matt Cramer",
    "Jennifer Mikar"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $names[$i]
}

# Fix the Id sequence for rows 100-181 (was a duplicated/garbled range in the source)
for ($row = 100; $row -le 181; $row++) {
    $ws.Cells.Item($row, 1).Value = 3104633 + ($row - 100)
}

# Move the active selection (was D27, scrolled to row 11)
$ws.Range("G3").Select() | Out-Null
